$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "16.79", "0.06556",
# "30.664.77"). Force those cells to Text format first so Excel stores the
# exact original string instead of re-parsing/rounding it as a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.664.77'
$ws.Range("E2").Value = '  +1.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.891.65'
$ws.Range("E3").Value = '  +2.00%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.63'
$ws.Range("E5").Value = '  +1.38%  '

$ws.Range("E7").Value = '  +1.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2868'
$ws.Range("E8").Value = '  +2.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06556'
$ws.Range("E9").Value = '  +1.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.920.39'
$ws.Range("E10").Value = '  +3.56%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.79'
$ws.Range("E11").Value = '  +3.55%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07465'
$ws.Range("E12").Value = '  +1.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.106'
$ws.Range("E13").Value = '  +0.36%  '

$ws.Range("E14").Value = '  +1.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6677'
$ws.Range("E15").Value = '  +3.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.656.17'
$ws.Range("E16").Value = '  +1.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.25'
$ws.Range("E17").Value = '  +1.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007576'
$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.163.18'
$ws.Range("E20").Value = '  +3.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '230.15'
$ws.Range("E21").Value = '  +1.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.272'
$ws.Range("E23").Value = '  -0.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.195'
$ws.Range("E24").Value = '  +1.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.07'
$ws.Range("E25").Value = '  +3.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.369'
$ws.Range("E26").Value = '  +1.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.87'
$ws.Range("E27").Value = '  +2.17%  '

$ws.Range("E28").Value = '  +2.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1026'
$ws.Range("E29").Value = '  +11.84%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.397'
$ws.Range("E30").Value = '  -2.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.337'
$ws.Range("E31").Value = '  +2.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.032'

$ws.Range("E33").Value = '  +2.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.213'
$ws.Range("E34").Value = '  +6.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7542'
$ws.Range("E35").Value = '  +3.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").Value = '  +0.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").Value = '  +0.99%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01877'
$ws.Range("E38").Value = '  +1.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.648'
$ws.Range("E39").Value = '  +2.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9191'
$ws.Range("E40").Value = '  +2.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.067'
$ws.Range("E41").Value = '  +0.97%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '107.05'
$ws.Range("E42").Value = '  +1.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4292'
$ws.Range("E43").Value = '  +1.47%  '

$ws.Range("E44").Value = '  +0.27%  '

$ws.Range("E45").Value = '  -4.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.429'
$ws.Range("E46").Value = '  +0.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '64.39'
$ws.Range("E47").Value = '  +0.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1276'
$ws.Range("E48").Value = '  -2.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.495'
$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.953'
$ws.Range("E50").Value = '  +2.63%  '

$ws.Range("E51").Value = '  +0.67%  '
